$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows of data for "Nome 1" and "Nome 2" (rows 2 and 3).
# This is a content clear (not a row delete/shift): rows 4-7 keep their
# original row numbers, matching the commit "retirada de duas instancias".
$ws.Range("A2:C3").ClearContents()

# Reflect the resulting selection left by the user's delete action.
$ws.Range("A2:XFD3").Select()
